$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binding_sites_numbering")

# New "7 (Pore)" binding-site rows appended after the existing data (rows 130-137).
# Columns: A=BS, B=GABRA1_Alphafold, C=6X40_D, D=5VDH_D, E=(highlight/blank)
$colA = @("7 (Pore)", "7 (Pore)", "7 (Pore)", "7 (Pore)", "7 (Pore)", "7 (Pore)", "7 (Pore)", "7 (Pore)")
$colB = @("P280", "V284", "T288", "L291", "T292", "T295", "I298", "S299")
$colC = @("P253", "V257", "T261", "L264", "T265", "T268", "I271", "S272")
$colD = @("P250", "A254", "T258", "L261", "T263", "T265", "S268", "G269")

$startRow = 130
$lastRow = $startRow + $colA.Length - 1

for ($row = $startRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 1).Value = $colA[$row - $startRow]
}
for ($row = $startRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 2).Value = $colB[$row - $startRow]
}
for ($row = $startRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 3).Value = $colC[$row - $startRow]
}
for ($row = $startRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 4).Value = $colD[$row - $startRow]
}

for ($row = $startRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 1).Font.Bold = $true
    $ws.Cells.Item($row, 1).VerticalAlignment = -4108

    $ws.Cells.Item($row, 2).HorizontalAlignment = -4108
    $ws.Cells.Item($row, 3).HorizontalAlignment = -4108
    $ws.Cells.Item($row, 4).HorizontalAlignment = -4108

    $ws.Cells.Item($row, 5).Interior.Color = 5296274
}

$ws.Range("I127").Select()
